$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update each coins Price (D) and 1h Volume
# change (E) cell to the latest snapshot. Both columns hold plain text in
# the source data (not numbers), e.g. "25.982.81" or "  +0.18%  ", so for
# any new Price value that Excel would otherwise auto-parse as a number we
# briefly force a text format before writing it, then clear the format again
# right away so the cell ends up with no explicit style -- same as every
# other untouched data cell in the sheet.

$ws.Range("D2").Value = '25.982.81'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.638.46'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.88'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5125'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2574'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07772'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").Value = '1.647.48'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5455'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '0.0₅7735'
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '25.997.79'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '196.52'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.417'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.903'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.084'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.924'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1228'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +6.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.840'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.239'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04844'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.270'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.202'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.535'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.376'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9129'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5539'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '1.099.72'
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01569'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.001'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.520'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.548'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8047'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.98'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("D45").Value = '0.0₈122'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").Value = '1.778.15'
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4538'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05211'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.498'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.02%  '
